$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (time slots shift from morning 07:30-12:00 to afternoon/evening 12:30-17:30) ---
$headers = @("12:30:00","13:00:00","13:30:00","14:00:00","14:30:00","15:00:00","15:30:00","16:00:00","16:30:00","17:00:00","17:30:00")
$cols = @("B","C","D","E","F","G","H","I","J","K","L")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $headers[$i]
}

# Copy header formatting (bold font, border, centered alignment) from K1 onto the new L1 header cell
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update room availability matrix for rows 2-44 (columns B-K) and add new column L ---
$avail = @{}
$avail[2] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[3] = @($false,$false,$false,$false,$false,$false,$false,$false,$false,$false,$false)
$avail[4] = @($true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false)
$avail[5] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[6] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[7] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[8] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[9] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[10] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[11] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[12] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[13] = @($true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false)
$avail[14] = @($true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false)
$avail[15] = @($true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false)
$avail[16] = @($true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false)
$avail[17] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[18] = @($true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false)
$avail[19] = @($true,$true,$true,$true,$false,$false,$false,$false,$false,$false,$false)
$avail[20] = @($true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false)
$avail[21] = @($true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false)
$avail[22] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[23] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[24] = @($true,$true,$true,$true,$false,$false,$false,$false,$false,$false,$false)
$avail[25] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[26] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[27] = @($true,$true,$true,$true,$false,$false,$false,$false,$false,$false,$false)
$avail[28] = @($true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false)
$avail[29] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[30] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[31] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[32] = @($true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false)
$avail[33] = @($true,$true,$true,$true,$false,$false,$false,$false,$false,$false,$false)
$avail[34] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[35] = @($true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false)
$avail[36] = @($true,$true,$true,$true,$false,$false,$false,$false,$false,$false,$false)
$avail[37] = @($true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false)
$avail[38] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[39] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[40] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[41] = @($true,$true,$true,$true,$false,$false,$false,$false,$false,$false,$false)
$avail[42] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[43] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)
$avail[44] = @($true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false)

foreach ($r in $avail.Keys) {
    $vals = $avail[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

Write-Host "Done updating schedule to afternoon/evening slots."
